# Generate Report for Handback
#
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# 5c709bbf-1c18-4af7-b3a8-efa64bb1f25d row (row 7) on both the
# zh-cn and de-de language sheets, now that a (stale) handback was
# received for that item.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c57b892698cf1c049db8b03eb9351d547c06b47/e2e/5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0229c306a1f85ef398d3ed6e1565a140ccad4de/e2e/5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md."

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("I7").Value = "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md"
$ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0229c306a1f85ef398d3ed6e1565a140ccad4de/e2e/5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md", "", "", "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md")

$ws.Range("J7").Value = "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.6764fd3ca2d5b07f86432096017c0188943ce414.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-18 08:57:28"
$ws.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("I7").Value = "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md"
$ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0229c306a1f85ef398d3ed6e1565a140ccad4de/e2e/5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md", "", "", "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.md")

$ws.Range("J7").Value = "5c709bbf-1c18-4af7-b3a8-efa64bb1f25d.6764fd3ca2d5b07f86432096017c0188943ce414.de-de.xlf"
$ws.Range("K7").Value = "2016-08-18 08:57:35"
$ws.Range("P7").Value = $errorDetail
